$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '60.816.73'
Set-TextValue "E2" '  +3.07%  '

# Row 3
Set-TextValue "D3" '2.713.33'
Set-TextValue "E3" '  +3.01%  '

# Row 4
Set-TextValue "D4" '0.999'
Set-TextValue "E4" '  -0.11%  '

# Row 5
Set-TextValue "D5" '526.58'
Set-TextValue "E5" '  +1.67%  '

# Row 6
Set-TextValue "D6" '144.81'
Set-TextValue "E6" '  +0.19%  '

# Row 7
Set-TextValue "E7" '  +0.12%  '

# Row 8
Set-TextValue "E8" '  +2.38%  '

# Row 9
Set-TextValue "D9" '2.713.97'
Set-TextValue "E9" '  +2.12%  '

# Row 10
Set-TextValue "D10" '6.62'
Set-TextValue "E10" '  +5.98%  '

# Row 11
Set-TextValue "E11" '  +1.36%  '

# Row 12
Set-TextValue "E12" '  +1.43%  '

# Row 13
Set-TextValue "E13" '  +2.85%  '

# Row 14
Set-TextValue "D14" '3.165.86'
Set-TextValue "E14" '  +2.11%  '

# Row 15
Set-TextValue "D15" '60.814.36'
Set-TextValue "E15" '  +3.17%  '

# Row 16
Set-TextValue "B16" 'Avalanche'
Set-TextValue "C16" 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D16" '21.31'
Set-TextValue "E16" '  +2.10%  '

# Row 17
Set-TextValue "B17" 'ShibaInu'
Set-TextValue "C17" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D17" '0.0000138'
Set-TextValue "E17" '  +0.96%  '

# Row 18
Set-TextValue "B18" 'WrappedEther'
Set-TextValue "C18" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D18" '2.710.51'
Set-TextValue "E18" '  +2.00%  '

# Row 19
Set-TextValue "D19" '348.85'
Set-TextValue "E19" '  +0.44%  '

# Row 20
Set-TextValue "E20" '  +0.03%  '

# Row 21
Set-TextValue "E21" '  +2.29%  '

# Row 22
Set-TextValue "D22" '6.35'
Set-TextValue "E22" '  +2.51%  '

# Row 23
Set-TextValue "E23" '  +0.31%  '

# Row 24
Set-TextValue "D24" '63.66'
Set-TextValue "E24" '  +3.10%  '

# Row 25
Set-TextValue "B25" 'Kaspa'
Set-TextValue "C25" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D25" '0.171'
Set-TextValue "E25" '  +6.38%  '

# Row 26
Set-TextValue "B26" 'Polygon'
Set-TextValue "C26" 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D26" '0.421'
Set-TextValue "E26" '  +0.53%  '

# Row 27
Set-TextValue "D27" '0.992'
Set-TextValue "E27" '  -0.12%  '

# Row 28
Set-TextValue "B28" 'PEPE'
Set-TextValue "C28" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D28" '0.0₃0821'
Set-TextValue "E28" '  +2.33%  '

# Row 29
Set-TextValue "B29" 'InternetComputer(DFINITY)'
Set-TextValue "C29" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D29" '7.32'
Set-TextValue "E29" '  +3.00%  '

# Row 30
Set-TextValue "E30" '  +9.37%  '

# Row 31
Set-TextValue "E31" '  +0.08%  '

# Row 32
Set-TextValue "E32" '  +1.91%  '

# Row 33
Set-TextValue "D33" '19.14'
Set-TextValue "E33" '  +1.07%  '

# Row 34
Set-TextValue "D34" '149.97'
Set-TextValue "E34" '  +0.28%  '

# Row 35
Set-TextValue "D35" '4.27'
Set-TextValue "E35" '  +6.49%  '

# Row 36
Set-TextValue "E36" '  +9.16%  '

# Row 37
Set-TextValue "D37" '0.946'
Set-TextValue "E37" '  -2.81%  '

# Row 38
Set-TextValue "E38" '  +4.94%  '

# Row 39
Set-TextValue "E39" '  +8.06%  '

# Row 40
Set-TextValue "D40" '36.99'
Set-TextValue "E40" '  +0.69%  '

# Row 41
Set-TextValue "D41" '3.68'
Set-TextValue "E41" '  -0.65%  '

# Row 42
Set-TextValue "D42" '287.24'
Set-TextValue "E42" '  +3.30%  '

# Row 43
Set-TextValue "D43" '0.0994'
Set-TextValue "E43" '  +1.13%  '

# Row 44
Set-TextValue "B44" 'EnergySwap'
Set-TextValue "C44" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D44" '20.03'
Set-TextValue "E44" '  +2.01%  '

# Row 45
Set-TextValue "B45" 'FirstDigitalUSD'
Set-TextValue "C45" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D45" '0.998'
Set-TextValue "E45" '  +0.41%  '

# Row 46
Set-TextValue "B46" 'Mantle'
Set-TextValue "C46" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D46" '0.612'
Set-TextValue "E46" '  +0.14%  '

# Row 47
Set-TextValue "D47" '2.141.43'
Set-TextValue "E47" '  +7.69%  '

# Row 48
Set-TextValue "D48" '0.0542'
Set-TextValue "E48" '  +2.26%  '

# Row 49
Set-TextValue "B49" 'VeChain'
Set-TextValue "C49" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D49" '0.0235'
Set-TextValue "E49" '  +2.65%  '

# Row 50
Set-TextValue "B50" 'RenderToken'
Set-TextValue "C50" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D50" '4.80'
Set-TextValue "E50" '  +2.76%  '

# Row 51
Set-TextValue "E51" '  +1.68%  '
